# Insert one new data row at row 59 (pushing the existing rows 59-150 down
# to 60-151), then populate the newly-inserted row 59 with its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(59).Insert()

$ws.Range("A59").Value = 9
$ws.Range("B59").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").Value = 44915
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = 100112022
$ws.Range("G59").Value = "Arveja Verde"
$ws.Range("H59").Value = "Perfection"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 34
$ws.Range("K59").Value = 18000
$ws.Range("L59").Value = 20000
$ws.Range("M59").Value = 19000
$ws.Range("N59").Value = '$/malla 25 kilos'
$ws.Range("O59").Value = "Carahue"
$ws.Range("P59").Value = 760
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"
